$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest scraped values.
# For D-column price cells, the raw text must be preserved exactly as scraped (it is
# sometimes not a valid Excel number, e.g. "62.366.98", and even when it is numeric-looking
# it must stay textual, e.g. "121.20" must not collapse to 121.2). We force text entry by
# pre-setting the cell number format to Text ("@"), assigning the value, then resetting the
# cell style back to "Normal" so no residual style/number-format is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.366.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.029.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.46%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.28%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.023.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.444"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000218"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.514.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("E16").Value = "  +1.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.328.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.024.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("E31").Value = "  +2.70%  "

$ws.Range("E32").Value = "  -2.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "58.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "465.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.196.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0386"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0780"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.28%  "

$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.42%  "

$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.245"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "121.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.23%  "

$ws.Range("E48").Value = "  +0.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0508"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("E51").Value = "  +4.60%  "
